$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.473.22"
$ws.Range("E2").Value = "  +1.61%  "
$ws.Range("D3").Value = "2.280.21"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D5").Value = "'307.40"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "'97.51"
$ws.Range("E6").Value = "  +5.29%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("D10").Value = "'36.02"
$ws.Range("E10").Value = "  +10.32%  "
$ws.Range("D11").Value = "'0.0797"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "2.634.94"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("D16").Value = "2.290.63"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "'0.798"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "42.378.11"
$ws.Range("E19").Value = "  +0.46%  "
$ws.Range("D20").Value = "0.0₃0911"
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'67.79"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'240.72"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").Value = "'1.95"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'23.91"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").Value = "'37.74"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("D29").Value = "'9.52"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").Value = "'159.18"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("D35").Value = "'0.0741"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("D36").Value = "'17.05"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "'0.105"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  +4.46%  "
$ws.Range("E42").Value = "  +14.09%  "
$ws.Range("D43").Value = "1.999.47"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("E45").Value = "  -2.81%  "
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "'9.99"
$ws.Range("E47").Value = "  -3.43%  "
$ws.Range("D48").Value = "'52.97"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'1.52"
$ws.Range("E49").Value = "  +0.52%  "
$ws.Range("D50").Value = "'72.19"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").Value = "'92.07"
$ws.Range("E51").Value = "  +0.87%  "
